# Generate Report for Handoff
# Updates the handoff timestamps (Overview, zh-cn, de-de sheets) and sets
# the "Priority" column to "ht" for the rows that were just handed off
# (everything except row 10, which corresponds to a different file/status).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 11, 12, 13)

foreach ($r in $rows) {
    # Overview sheet: "Latest HO Xliff Generate Date" (column G)
    $wsOverview.Range("G$r").Value = "2016-08-31 04:21:57"

    # zh-cn sheet: "Priority" (column E) and "Latest Handoff Datetime" (column H)
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-31 04:21:50"

    # de-de sheet: "Priority" (column E) and "Latest Handoff Datetime" (column H)
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-31 04:21:57"
}
